$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-09-23"

# Update the September label in column A (row 10) to reflect the new date
$ws.Range("A10").Value = "September (through 09-23)"

# August row (row 9) - only 2021 (H) column changes
$ws.Range("H9").Value = 159

# September row (row 10) - update 2017-2021 values (2016/C10 unchanged)
$ws.Range("B10").Value = 25
$ws.Range("D10").Value = 57
$ws.Range("E10").Value = 46
$ws.Range("F10").Value = 57
$ws.Range("G10").Value = 87
$ws.Range("H10").Value = 135

# Total row (row 11) - update 2017-2021 values (2016/C11 unchanged)
$ws.Range("B11").Value = 219
$ws.Range("D11").Value = 608
$ws.Range("E11").Value = 536
$ws.Range("F11").Value = 406
$ws.Range("G11").Value = 871
$ws.Range("H11").Value = 1206
